$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "place name" column (B). This shifts the salary values
# that used to live in column C left into column B.
$ws.Columns("B").Delete()

# Fix up the style of A33 (it previously relied on a now-unused style,
# so line it up with the plain style used by its neighbours on that row).
$ws.Range("B33").Copy()
$ws.Range("A33").PasteSpecial(-4122)

# Reflect the new, narrower data range in the selection and print area.
$ws.Columns("B").EntireColumn.Select()
$ws.PageSetup.PrintArea = '$B$1:$B$21'
